$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-07 01:59:46"

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
